# Add PF/1.0.3 to meta-sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "PF/1.0.3"
$ws.Range("B3").Value = "X"
$ws.Range("C3").Value = "X"
$ws.Range("D3").Value = "X"

# New row has no explicit style in the target (unlike the styled header /
# value rows above it), so reset it back to the workbook default style
# instead of inheriting the column's style.
$ws.Range("A3:D3").Style = "Normal"
